$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each data row (2-20), shift existing values in columns B..K one column to the
# right (dropping whatever falls out past column K), then write the new forecast-error
# value into column B. Columns are updated from right to left so a value is always
# read before it gets overwritten.

$ws.Range("K2").Value = $ws.Range("J2").Value2
$ws.Range("J2").Value = $ws.Range("I2").Value2
$ws.Range("I2").Value = $ws.Range("H2").Value2
$ws.Range("H2").Value = $ws.Range("G2").Value2
$ws.Range("G2").Value = $ws.Range("F2").Value2
$ws.Range("F2").Value = $ws.Range("E2").Value2
$ws.Range("E2").Value = $ws.Range("D2").Value2
$ws.Range("D2").Value = $ws.Range("C2").Value2
$ws.Range("C2").Value = $ws.Range("B2").Value2
$ws.Range("B2").Value = 0.3719860057927588

$ws.Range("K3").Value = $ws.Range("J3").Value2
$ws.Range("J3").Value = $ws.Range("I3").Value2
$ws.Range("I3").Value = $ws.Range("H3").Value2
$ws.Range("H3").Value = $ws.Range("G3").Value2
$ws.Range("G3").Value = $ws.Range("F3").Value2
$ws.Range("F3").Value = $ws.Range("E3").Value2
$ws.Range("E3").Value = $ws.Range("D3").Value2
$ws.Range("D3").Value = $ws.Range("C3").Value2
$ws.Range("C3").Value = $ws.Range("B3").Value2
$ws.Range("B3").Value = -2.702915518772638

$ws.Range("K4").Value = $ws.Range("J4").Value2
$ws.Range("J4").Value = $ws.Range("I4").Value2
$ws.Range("I4").Value = $ws.Range("H4").Value2
$ws.Range("H4").Value = $ws.Range("G4").Value2
$ws.Range("G4").Value = $ws.Range("F4").Value2
$ws.Range("F4").Value = $ws.Range("E4").Value2
$ws.Range("E4").Value = $ws.Range("D4").Value2
$ws.Range("D4").Value = $ws.Range("C4").Value2
$ws.Range("C4").Value = $ws.Range("B4").Value2
$ws.Range("B4").Value = -0.2307826431404359

$ws.Range("K5").Value = $ws.Range("J5").Value2
$ws.Range("J5").Value = $ws.Range("I5").Value2
$ws.Range("I5").Value = $ws.Range("H5").Value2
$ws.Range("H5").Value = $ws.Range("G5").Value2
$ws.Range("G5").Value = $ws.Range("F5").Value2
$ws.Range("F5").Value = $ws.Range("E5").Value2
$ws.Range("E5").Value = $ws.Range("D5").Value2
$ws.Range("D5").Value = $ws.Range("C5").Value2
$ws.Range("C5").Value = $ws.Range("B5").Value2
$ws.Range("B5").Value = -0.5654386276933741

$ws.Range("K6").Value = $ws.Range("J6").Value2
$ws.Range("J6").Value = $ws.Range("I6").Value2
$ws.Range("I6").Value = $ws.Range("H6").Value2
$ws.Range("H6").Value = $ws.Range("G6").Value2
$ws.Range("G6").Value = $ws.Range("F6").Value2
$ws.Range("F6").Value = $ws.Range("E6").Value2
$ws.Range("E6").Value = $ws.Range("D6").Value2
$ws.Range("D6").Value = $ws.Range("C6").Value2
$ws.Range("C6").Value = $ws.Range("B6").Value2
$ws.Range("B6").Value = -0.6603092772102132

$ws.Range("K7").Value = $ws.Range("J7").Value2
$ws.Range("J7").Value = $ws.Range("I7").Value2
$ws.Range("I7").Value = $ws.Range("H7").Value2
$ws.Range("H7").Value = $ws.Range("G7").Value2
$ws.Range("G7").Value = $ws.Range("F7").Value2
$ws.Range("F7").Value = $ws.Range("E7").Value2
$ws.Range("E7").Value = $ws.Range("D7").Value2
$ws.Range("D7").Value = $ws.Range("C7").Value2
$ws.Range("C7").Value = $ws.Range("B7").Value2
$ws.Range("B7").Value = -0.15162438770796

$ws.Range("K8").Value = $ws.Range("J8").Value2
$ws.Range("J8").Value = $ws.Range("I8").Value2
$ws.Range("I8").Value = $ws.Range("H8").Value2
$ws.Range("H8").Value = $ws.Range("G8").Value2
$ws.Range("G8").Value = $ws.Range("F8").Value2
$ws.Range("F8").Value = $ws.Range("E8").Value2
$ws.Range("E8").Value = $ws.Range("D8").Value2
$ws.Range("D8").Value = $ws.Range("C8").Value2
$ws.Range("C8").Value = $ws.Range("B8").Value2
$ws.Range("B8").Value = -0.2053460154962278

$ws.Range("K9").Value = $ws.Range("J9").Value2
$ws.Range("J9").Value = $ws.Range("I9").Value2
$ws.Range("I9").Value = $ws.Range("H9").Value2
$ws.Range("H9").Value = $ws.Range("G9").Value2
$ws.Range("G9").Value = $ws.Range("F9").Value2
$ws.Range("F9").Value = $ws.Range("E9").Value2
$ws.Range("E9").Value = $ws.Range("D9").Value2
$ws.Range("D9").Value = $ws.Range("C9").Value2
$ws.Range("C9").Value = $ws.Range("B9").Value2
$ws.Range("B9").Value = 0.6162032393936197

$ws.Range("K10").Value = $ws.Range("J10").Value2
$ws.Range("J10").Value = $ws.Range("I10").Value2
$ws.Range("I10").Value = $ws.Range("H10").Value2
$ws.Range("H10").Value = $ws.Range("G10").Value2
$ws.Range("G10").Value = $ws.Range("F10").Value2
$ws.Range("F10").Value = $ws.Range("E10").Value2
$ws.Range("E10").Value = $ws.Range("D10").Value2
$ws.Range("D10").Value = $ws.Range("C10").Value2
$ws.Range("C10").Value = $ws.Range("B10").Value2
$ws.Range("B10").Value = 1.652643173475852

$ws.Range("K11").Value = $ws.Range("J11").Value2
$ws.Range("J11").Value = $ws.Range("I11").Value2
$ws.Range("I11").Value = $ws.Range("H11").Value2
$ws.Range("H11").Value = $ws.Range("G11").Value2
$ws.Range("G11").Value = $ws.Range("F11").Value2
$ws.Range("F11").Value = $ws.Range("E11").Value2
$ws.Range("E11").Value = $ws.Range("D11").Value2
$ws.Range("D11").Value = $ws.Range("C11").Value2
$ws.Range("C11").Value = $ws.Range("B11").Value2
$ws.Range("B11").Value = 0.3110387314724781

$ws.Range("J12").Value = $ws.Range("I12").Value2
$ws.Range("I12").Value = $ws.Range("H12").Value2
$ws.Range("H12").Value = $ws.Range("G12").Value2
$ws.Range("G12").Value = $ws.Range("F12").Value2
$ws.Range("F12").Value = $ws.Range("E12").Value2
$ws.Range("E12").Value = $ws.Range("D12").Value2
$ws.Range("D12").Value = $ws.Range("C12").Value2
$ws.Range("C12").Value = $ws.Range("B12").Value2
$ws.Range("B12").Value = 0.2388379152847414

$ws.Range("I13").Value = $ws.Range("H13").Value2
$ws.Range("H13").Value = $ws.Range("G13").Value2
$ws.Range("G13").Value = $ws.Range("F13").Value2
$ws.Range("F13").Value = $ws.Range("E13").Value2
$ws.Range("E13").Value = $ws.Range("D13").Value2
$ws.Range("D13").Value = $ws.Range("C13").Value2
$ws.Range("C13").Value = $ws.Range("B13").Value2
$ws.Range("B13").Value = 0.6508000635779043

$ws.Range("H14").Value = $ws.Range("G14").Value2
$ws.Range("G14").Value = $ws.Range("F14").Value2
$ws.Range("F14").Value = $ws.Range("E14").Value2
$ws.Range("E14").Value = $ws.Range("D14").Value2
$ws.Range("D14").Value = $ws.Range("C14").Value2
$ws.Range("C14").Value = $ws.Range("B14").Value2
$ws.Range("B14").Value = 0.2387740594105157

$ws.Range("G15").Value = $ws.Range("F15").Value2
$ws.Range("F15").Value = $ws.Range("E15").Value2
$ws.Range("E15").Value = $ws.Range("D15").Value2
$ws.Range("D15").Value = $ws.Range("C15").Value2
$ws.Range("C15").Value = $ws.Range("B15").Value2
$ws.Range("B15").Value = 0.3465902496671606

$ws.Range("F16").Value = $ws.Range("E16").Value2
$ws.Range("E16").Value = $ws.Range("D16").Value2
$ws.Range("D16").Value = $ws.Range("C16").Value2
$ws.Range("C16").Value = $ws.Range("B16").Value2
$ws.Range("B16").Value = 0.00230005330798793

$ws.Range("E17").Value = $ws.Range("D17").Value2
$ws.Range("D17").Value = $ws.Range("C17").Value2
$ws.Range("C17").Value = $ws.Range("B17").Value2
$ws.Range("B17").Value = -0.1902738424076751

$ws.Range("D18").Value = $ws.Range("C18").Value2
$ws.Range("C18").Value = $ws.Range("B18").Value2
$ws.Range("B18").Value = -0.3325070745318338

$ws.Range("C19").Value = $ws.Range("B19").Value2
$ws.Range("B19").Value = 0.1656141382254278

$ws.Range("B20").Value = -0.09587373626955231
Write-Output "done"
